$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (columns E,F,G,H,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.28161
$ws.Range("H2").Value = 0.84483
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4857986666666667
$ws.Range("N2").Value = 1.457396
$ws.Range("O2").Value = 0.4938122518903432
$ws.Range("P2").Value = 0.4938122518903432
$ws.Range("Q2").Value = 0.13680576252
$ws.Range("R2").Value = 1.23125186268
$ws.Range("S2").Value = 0.4938122518903432
$ws.Range("T2").Value = 0.4938122518903432

# Update existing row 3 values
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.28161
$ws.Range("H3").Value = 0.84483
$ws.Range("M3").Value = 0.4477856666666667
$ws.Range("O3").Value = 0.4551722011468782
$ws.Range("P3").Value = 0.4551722011468782
$ws.Range("Q3").Value = 0.12610092159
$ws.Range("R3").Value = 1.13490829431
$ws.Range("S3").Value = 0.4551722011468782
$ws.Range("T3").Value = 0.4551722011468782

# Add new row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Adm2"
$ws.Range("C4").Value = "Ramp3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.28161
$ws.Range("H4").Value = 0.84483
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05018766666666667
$ws.Range("N4").Value = 0.150563
$ws.Range("O4").Value = 0.05101554696277864
$ws.Range("P4").Value = 0.05101554696277864
$ws.Range("Q4").Value = 0.01413334881
$ws.Range("R4").Value = 0.12720013929
$ws.Range("S4").Value = 0.05101554696277864
$ws.Range("T4").Value = 0.05101554696277864
